$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID Resources-HCP")

# --- Fix the auto-flow height of two existing rows (113 and 123) ---
# These rows were re-flowed to a taller height in the authored edit.
$ws.Rows.Item(113).RowHeight = 45
$ws.Rows.Item(123).RowHeight = 45

# --- New row 126: BC Health Care Assistants needs-assessment report ---
$ws.Range("A126").Value = "All"
$ws.Range("B126").Value = "Healthcare Provider Wellness"
$ws.Range("C126").Value = "British Columbia"
$ws.Range("C126").Style = "Normal"
$ws.Range("C126").HorizontalAlignment = -4131
$ws.Range("D126").Value = "B.C. Health Care Assistants:  Assessment of education and support needs during COVID-19 Report- June 5, 2020"
$ws.Range("E126").Value = "BC Centre for Palliative Care and Life and Death Matters"
$ws.Range("E126").HorizontalAlignment = -4131
$ws.Range("F126").Value = "Publication"
$ws.Range("F126").HorizontalAlignment = -4131
$ws.Range("G126").Value = "https://bc-cpc.ca/cpc/wp-content/uploads/2020/03/HCA-needs-assessment-June-2020.pdf"
$ws.Hyperlinks.Add($ws.Range("G126"), "https://bc-cpc.ca/cpc/wp-content/uploads/2020/03/HCA-needs-assessment-June-2020.pdf")
$ws.Rows.Item(126).RowHeight = 45

# --- New row 127: Canadian Virtual Hospice grief module for healthcare workers ---
$ws.Range("A127").Value = "All"
$ws.Range("B127").Value = "Psychosocial Care"
$ws.Range("C127").Value = "Pan-Canadian"
$ws.Range("C127").HorizontalAlignment = -4131
$ws.Range("D127").Value = "For people working in healthcare:  Module 1 - COVID-19 and grief`n"
$ws.Range("E127").Value = "Canadian Virtual Hospice"
$ws.Range("E127").HorizontalAlignment = -4131
$ws.Range("F127").Value = "Education"
$ws.Range("F127").HorizontalAlignment = -4131
$link127 = "https://mygrief.azurewebsites.net/mod/lesson/view.php?id=554&utm_source=Canadian+Virtual+Hospice+%7C+Portail+canadien+en+soins+palliatifs&utm_campaign=9d332db16d-EMAIL_CAMPAIGN_2020_Oct30MonthlyEnews_COPY_01&utm_medium=email&utm_term=0_8030ddb8bb-9d332db16d-427597322"
$ws.Range("G127").Value = $link127
$ws.Hyperlinks.Add($ws.Range("G127"), $link127)
$ws.Range("G127").Style = "Hyperlink"
$ws.Rows.Item(127).RowHeight = 45
